# Applies the BPI-offer queries docx edit described by the commit
# "Adding more benchmark domains":
#   - fixes/typo edits on two bullet lines
#   - de-splits several query bullet lines whose runs/proofErr wrappers
#     had fragmented a single sentence into multiple runs
#   - re-homes the _GoBack bookmark that sat mid-sentence in the
#     "retrieve applications..." bullet to the start of the following
#     bullet ("which offer IDs...")

$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $range = $d.Content
    $ok = $range.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
    if (-not $ok) {
        Write-Output "WARNING: not found -> $find"
    }
    return $ok
}

# 1) "Tell" / " me" / " the rejected application..." -> single run, drop "me"
Replace-Text "Tell me the rejected application which has more than 1000 withdraw loan amount." "Tell the rejected application which has more than 1000 withdraw loan amount."

# 2) "...monthly cost larger than 100 since 2019" -> "...monthly coast larger than 100 since 2019"
Replace-Text "Show the approved application with monthly cost larger than 100 since 2019" "Show the approved application with monthly coast larger than 100 since 2019"

# 3) "list all " + "offer" + " ids with offered amount more than 10000" -> single run (drop proofErr wrap)
Replace-Text "list all offer ids with offered amount more than 10000" "list all offer ids with offered amount more than 10000"

# 4) "which applications had offered amount less than 1000 but credit score more than " + "1000" -> single run
Replace-Text "which applications had offered amount less than 1000 but credit score more than 1000" "which applications had offered amount less than 1000 but credit score more than 1000"

# 5) "What is the average offered amount which are accepted in last quarter of " + "2016" -> single run
Replace-Text "What is the average offered amount which are accepted in last quarter of 2016" "What is the average offered amount which are accepted in last quarter of 2016"

# 6) "what is the maximum monthly cost for offers completed before June " + "2016" -> single run
Replace-Text "what is the maximum monthly cost for offers completed before June 2016" "what is the maximum monthly cost for offers completed before June 2016"

# 7) "what is the minimum credit score across " + "offers" -> single run
Replace-Text "what is the minimum credit score across offers" "what is the minimum credit score across offers"

# 8) "retri" + "e" + "ve appli" + [bookmark] + "cations with maximum monthly cost more than 1000 in 2016"
#    -> "retri" + "e" + "ve applications with maximum monthly cost more than 1000 in 2016"
#    (Find/Replace treats the bookmark-interrupted text as contiguous and
#     merges the trailing two runs into one, dropping the now-redundant
#     bookmark from this mid-sentence spot.)
Replace-Text "ve applications with maximum monthly cost more than 1000 in 2016" "ve applications with maximum monthly cost more than 1000 in 2016"

# Re-home the _GoBack bookmark to the very start of the next bullet
# ("which offer IDs had average offered amount of less than 10000 in Q3, 2016").
$target = $d.Content
$found = $target.Find.Execute("which offer IDs had average offered amount of less than 10000 in Q3, 2016")
if ($found) {
    $target.Collapse(1)
    $d.Bookmarks.Add("_GoBack", $target)
} else {
    Write-Output "WARNING: could not find bookmark target paragraph"
}

Write-Output "done"
